$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 136-137, shifting the existing rows (old 136..151)
# down to 138..153.
$ws.Range("A136:R137").Insert()

# Populate the newly inserted rows with the new weekly record
# (Coliflor, Primera/Segunda) for date 44476 (2021-10-07).
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44476
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112008
$ws.Range("G136").Value = "Coliflor"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 600
$ws.Range("L136").Value = 700
$ws.Range("M136").Value = 650
$ws.Range("N136").Value = "$/unidad"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 650
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"

$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44476
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112008
$ws.Range("G137").Value = "Coliflor"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 500
$ws.Range("L137").Value = 500
$ws.Range("M137").Value = 500
$ws.Range("N137").Value = "$/unidad"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 500
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"
